# Update workbook/sheet title and shared header text to reflect data
# refreshed through December 08, 2021 (was December 07).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Through 2021-12-08"
$ws.Range("B1").Value = "December 2021 (through December 08)"

# Update/insert carjacking counts by neighborhood (row) / month (column)
# reflecting the additional day of reported incidents.

# West Town
$ws.Range("AX2").Value = 2     # December 2017: 1 -> 2

# Englewood
$ws.Range("B3").Value = 4      # December 2021 (through Dec 08): 2 -> 4
$ws.Range("BV3").Value = 1     # December 2015: new -> 1

# North Lawndale
$ws.Range("N4").Value = 6      # December 2020: 5 -> 6
$ws.Range("Z4").Value = 1      # December 2019: new -> 1

# Austin
$ws.Range("N7").Value = 4      # December 2020: 2 -> 4

# Chatham
$ws.Range("N8").Value = 1      # December 2020: new -> 1

# Grand Crossing
$ws.Range("B9").Value = 4      # December 2021 (through Dec 08): 3 -> 4

# Humboldt Park
$ws.Range("AX11").Value = 3    # December 2017: 2 -> 3
$ws.Range("BV11").Value = 1    # December 2015: new -> 1

# Roseland
$ws.Range("BJ13").Value = 3    # December 2016: 2 -> 3

# Wicker Park
$ws.Range("B21").Value = 1     # December 2021 (through Dec 08): new -> 1

# Kenwood
$ws.Range("AX22").Value = 2    # December 2017: 1 -> 2

# Little Village
$ws.Range("AL23").Value = 1    # December 2018: new -> 1

# South Shore
$ws.Range("AX24").Value = 1    # December 2017: new -> 1

# Lake View
$ws.Range("AL26").Value = 1    # December 2018: new -> 1

# Near South Side
$ws.Range("B33").Value = 2     # December 2021 (through Dec 08): 1 -> 2

# Woodlawn
$ws.Range("B34").Value = 3     # December 2021 (through Dec 08): 2 -> 3

# Auburn Gresham
$ws.Range("AL38").Value = 1    # December 2018: new -> 1
$ws.Range("AX38").Value = 4    # December 2017: 3 -> 4

# East Village
$ws.Range("B42").Value = 2     # December 2021 (through Dec 08): 1 -> 2

# Logan Square
$ws.Range("B82").Value = 2     # December 2021 (through Dec 08): 1 -> 2

# O'Hare
$ws.Range("Z91").Value = 1     # December 2019: new -> 1
